$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Grades")
$lab = $wb.Worksheets.Item("Lab-Roster")

# --- Column P/Q header ---
$ws.Range("P1").Value = "Brief Submission"

# --- Row 9 (Bush, Edward Cole) — late submission note + citation ---
$ws.Range("P9").Value = "Late 2 hours"
$ws.Range("Q9").Value = "Towards a balanced presentation and objective interpretation of acoustic and trawl survey data, with specific reference to the eastern Scotian Shelf"

# --- Row 23 (Waltman, Billie Nicole) ---
$ws.Range("P23").Value = "x"
$ws.Range("Q23").Value = "The effects of fisheries management practises on freshwater ecosystems"

# --- Row 12 (Mills, Logan Chandler) ---
$ws.Range("P12").Value = "x"
$ws.Range("Q12").Value = "Using reverse-time egg transport analysis for predicting Asian carp spawning grounds in the Illinois River"

# --- Row 3 (Bacon, Corey Lonell) ---
$ws.Range("P3").Value = "x"
$ws.Range("Q3").Value = "Implications of piscine predator control on the federally listed fountain darter. "

# --- Row 11 (Lampert, Shaina Shaira) ---
$ws.Range("P11").Value = "x"
$ws.Range("Q11").Value = "Management issues in the Lake Victoria watershed"

# --- Row 6 (Dziamniski, Matthew Joseph) ---
$ws.Range("P6").Value = "x"
$ws.Range("Q6").Value = "Growth response of largemouth bass (Micropterus salmoides) to catch-and-release angling: a 27-year mark–recapture study"

# --- Row 20 (Stolz, Emily Carole) ---
$ws.Range("P20").Value = "x"
$ws.Range("Q20").Value = "Tracking bowfin with acoustic telemetry: Insight into the ecology of a living fossil"

# --- Row 10 (Kitaif, Jude Carson) ---
$ws.Range("P10").Value = "x"
$ws.Range("Q10").Value = "Contrasting patterns of productivity and survival rates for stream-type chinook salmon (Oncorhynchus tshawytscha) populations of the Snake and Columbia rivers"

# --- Row 4 (Burran, Sara Gabrielle) ---
$ws.Range("P4").Value = "x"
$ws.Range("Q4").Value = "Relatedness and body size influence territorial behaviour in Salmo salar juveniles in the wild."

# --- Row 8 (Holiman, Laura Haley) ---
$ws.Range("P8").Value = "x"
$ws.Range("Q8").Value = "Red Snapper Distribution on Natural Habitats and Artificial Structures in the Northern Gulf of Mexico"

# --- Row 14 (Norman, Durham Alexander) ---
$ws.Range("P14").Value = "x"
$ws.Range("Q14").Value = "Temperature and hydrologic alteration predict the spread of invasive Largemouth Bass (Micropterus salmoides)"

# --- Row 22 (Taylor, Kenneth Brandon) ---
$ws.Range("P22").Value = "x"
$ws.Range("Q22").Value = "Predictive Evaluation of Size Restrictions as Management Strategies for Tennessee Reservoir Crappie Fisheries"

# --- Row 7 (Godbey, Carice Nicole) ---
$ws.Range("P7").Value = "x"
$ws.Range("Q7").Value = "Accounting for variable recruitment and fishing mortality in 1 length-based stock assessments for data-limited fisheries"

# --- Row 16 (Rodgers, Colton Shane) ---
$ws.Range("P16").Value = "x"
$ws.Range("Q16").Value = "GIS visualisation and analysis of mobile hydroacoustic fisheries data: a practical example"

# --- Row 2 (Amacker, Caleb Agustus) ---
$ws.Range("P2").Value = "x"
$ws.Range("Q2").Value = "Public Perception of Agricultural Pollution and Gulf of Mexico Hypoxia"

# --- Row 21 (Taylor, Daniel Roane) ---
$ws.Range("P21").Value = "x"
$ws.Range("Q21").Value = "Assessing a social norms approach for improving recreational fisheries compliance"

# --- Row 15 (Red, Paige Delane) ---
$ws.Range("P15").Value = "x"
$ws.Range("Q15").Value = "Effects of hot dry summers on the loss of Atlantic salmon, Salmo salar, from estuaries in South West England"

# --- Row 13 (Moore, Jacob Andrew) ---
$ws.Range("P13").Value = "x"
$ws.Range("Q13").Value = "Comparing commercial and recreational harvest characteristics of paddlefish Polyodon spathula (Walbaum, 1792) in the Middle Mississippi River"

# --- Row 24 (Winterhalter, Emma Kiley) ---
$ws.Range("P24").Value = "x"
$ws.Range("Q24").Value = "Effects of Multiple Low-Head Dams on Fish, Macroinvertebrates, Habitat, and Water Quality in the Fox River, Illinois"

# --- Row 18 (Smith, Ashton Laray) ---
$ws.Range("P18").Value = "x"
$ws.Range("Q18").Value = "Effectively managing angler satisfaction in recreational fisheries requires understanding the fish species and the anglers"

# --- Row 5 (Bush, Edward Cole's neighbour — Dziamniski row above was 6; row 5 is Bush) ---
$ws.Range("P5").Value = "xx"

# --- Row 25 (Norris, David M) ---
$ws.Range("P25").Value = "x"
$ws.Range("Q25").Value = "Influence of behavior and mating success on brood-specific contribution to fish recruitment in ponds"

# --- Row 19 (Smith, Jonathan Charles) — numeric 0, not text ---
$ws.Range("P19").Value = 0

# --- Sheet-scoped defined name "citation" on the Grades sheet ---
$ws.Names.Add("citation", "=Grades!`$Q`$3")

# --- Page setup tweak on Grades sheet ---
$ws.PageSetup.Orientation = 1

# --- Selection state ---
$ws.Activate()
$ws.Range("P22").Select()

# --- Lab-Roster sheet: column widths + selection ---
$lab.Columns.Item(5).ColumnWidth = 8.333333333333332
$lab.Columns.Item(6).ColumnWidth = 7.666666666666667
$lab.Activate()
$lab.Range("H5").Select()

# restore Grades as the active sheet/tab (matches tabSelected on Grades in target)
$ws.Activate()
